$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("功能表")

# Update B7 value from "X" to "O" (feature now implemented:
# checking that a newly selected course does not share a name
# with an already-selected course)
$ws.Range("B7").Value = "O"

# Move the active selection from B8 to B9
$ws.Range("B9").Select()
